# Update row 20 (2025Q2) metrics on Sheet1 to reflect latest recorrencia data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C20").Value = 277
$ws.Range("D20").Value = 228
$ws.Range("E20").Value = 49
$ws.Range("F20").Value = 75
